$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F holding "value / 1000" helper figures next to the two
# measurement blocks (row 4 = "cpu" row of the first table, row 14 =
# "cpu" row of the third table). Formatted in scientific notation.
$ws.Range("F4").Formula = "=B4/1000"
$ws.Range("F4").NumberFormat = "0.00E+00"

$ws.Range("F14").Formula = "=B14/1000"
$ws.Range("F14").NumberFormat = "0.00E+00"

# Give the new column roughly the same width the author dragged it to.
$ws.Columns("F").ColumnWidth = 12.3

# Drop now-pointless empty placeholder cells that used to pad out rows
# 11, 15 and 16 (A11/D11, all of row 15, A16/B16/D16).
$ws.Range("A11").ClearContents()
$ws.Range("D11").ClearContents()

$ws.Range("A15:D15").ClearContents()

$ws.Range("A16").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("D16").ClearContents()
